# Scheduled-runner style data refresh for the Kujata_Profits workbook.
# Updates computed market-price / profit columns (H-N) for a set of
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# mirroring the values produced by the periodic price-scraper run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Cells.Item(112, 8).Value2 = 2584.6072
$ws.Cells.Item(112, 10).Value2 = 2698.8462
$ws.Cells.Item(112, 12).Value2 = 8096.5386
$ws.Cells.Item(112, 14).Value2 = -10312.5386

# Row 113: Amaro Kart / Starch Glue
$ws.Cells.Item(113, 8).Value2 = 20003152
$ws.Cells.Item(113, 9).Value2 = 40002744
$ws.Cells.Item(113, 10).Value2 = 3562
$ws.Cells.Item(113, 11).Value2 = 40002744
$ws.Cells.Item(113, 12).Value2 = 3562
$ws.Cells.Item(113, 13).Value2 = -39999490
$ws.Cells.Item(113, 14).Value2 = -10070

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value2 = 1613.625
$ws.Cells.Item(137, 9).Value2 = 1330.0476
$ws.Cells.Item(137, 10).Value2 = 2155
$ws.Cells.Item(137, 11).Value2 = 3990.142800000001
$ws.Cells.Item(137, 12).Value2 = 6465
$ws.Cells.Item(137, 13).Value2 = -1440.142800000001
$ws.Cells.Item(137, 14).Value2 = -11565

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value2 = 1486.8889
$ws.Cells.Item(141, 9).Value2 = 1486.8889
$ws.Cells.Item(141, 10).Value2 = 0
$ws.Cells.Item(141, 11).Value2 = 4460.6667
$ws.Cells.Item(141, 12).Value2 = 0
$ws.Cells.Item(141, 13).Value2 = 719.3333000000002
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value2 = 5210.087
$ws.Cells.Item(2, 9).Value2 = 701.4545000000001
$ws.Cells.Item(2, 10).Value2 = 9343
$ws.Cells.Item(2, 11).Value2 = 701.4545000000001
$ws.Cells.Item(2, 12).Value2 = 9343
$ws.Cells.Item(2, 13).Value2 = -588.4545000000001
$ws.Cells.Item(2, 14).Value2 = -9569

# Row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value2 = 4109.385
$ws.Cells.Item(32, 9).Value2 = 4169.892
$ws.Cells.Item(32, 11).Value2 = 4169.892
$ws.Cells.Item(32, 13).Value2 = -3882.892

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Cells.Item(45, 8).Value2 = 2246.2
$ws.Cells.Item(45, 9).Value2 = 2246.2
$ws.Cells.Item(45, 10).Value2 = 0
$ws.Cells.Item(45, 11).Value2 = 2246.2
$ws.Cells.Item(45, 12).Value2 = 0
$ws.Cells.Item(45, 13).Value2 = -1869.2
$ws.Cells.Item(45, 14).ClearContents()

# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Cells.Item(46, 8).Value2 = 5047.75
$ws.Cells.Item(46, 10).Value2 = 3400
$ws.Cells.Item(46, 12).Value2 = 3400
$ws.Cells.Item(46, 14).Value2 = -4038

# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value2 = 5210.087
$ws.Cells.Item(116, 9).Value2 = 701.4545000000001
$ws.Cells.Item(116, 10).Value2 = 9343
$ws.Cells.Item(116, 11).Value2 = 701.4545000000001
$ws.Cells.Item(116, 12).Value2 = 9343
$ws.Cells.Item(116, 13).Value2 = 1592.5455
$ws.Cells.Item(116, 14).Value2 = -13931

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value2 = 5210.087
$ws.Cells.Item(3, 9).Value2 = 701.4545000000001
$ws.Cells.Item(3, 10).Value2 = 9343
$ws.Cells.Item(3, 11).Value2 = 701.4545000000001
$ws.Cells.Item(3, 12).Value2 = 9343
$ws.Cells.Item(3, 13).Value2 = -587.4545000000001
$ws.Cells.Item(3, 14).Value2 = -9571

# Row 6: The Unkindest Cut / Bronze Saw
$ws.Cells.Item(6, 8).Value2 = 22500
$ws.Cells.Item(6, 10).Value2 = 22500
$ws.Cells.Item(6, 12).Value2 = 22500
$ws.Cells.Item(6, 14).Value2 = -22726

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Cells.Item(99, 8).Value2 = 33334348
$ws.Cells.Item(99, 9).Value2 = 33334348
$ws.Cells.Item(99, 11).Value2 = 33334348
$ws.Cells.Item(99, 13).Value2 = -33332850

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value2 = 2428
$ws.Cells.Item(107, 9).Value2 = 1803.6666
$ws.Cells.Item(107, 10).Value2 = 2802.6
$ws.Cells.Item(107, 11).Value2 = 1803.6666
$ws.Cells.Item(107, 12).Value2 = 2802.6
$ws.Cells.Item(107, 13).Value2 = 116.3334
$ws.Cells.Item(107, 14).Value2 = -6642.6

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value2 = 1411.9762
$ws.Cells.Item(31, 9).Value2 = 1048.091
$ws.Cells.Item(31, 10).Value2 = 1812.25
$ws.Cells.Item(31, 11).Value2 = 1048.091
$ws.Cells.Item(31, 12).Value2 = 1812.25
$ws.Cells.Item(31, 13).Value2 = -753.0909999999999
$ws.Cells.Item(31, 14).Value2 = -2402.25

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value2 = 1411.9762
$ws.Cells.Item(34, 9).Value2 = 1048.091
$ws.Cells.Item(34, 10).Value2 = 1812.25
$ws.Cells.Item(34, 11).Value2 = 1048.091
$ws.Cells.Item(34, 12).Value2 = 1812.25
$ws.Cells.Item(34, 13).Value2 = -846.0909999999999
$ws.Cells.Item(34, 14).Value2 = -2216.25

# Row 99: O Pine / Pine Lumber
$ws.Cells.Item(99, 8).Value2 = 1463482.6
$ws.Cells.Item(99, 9).Value2 = 1881083.8
$ws.Cells.Item(99, 10).Value2 = 1878.5
$ws.Cells.Item(99, 11).Value2 = 1881083.8
$ws.Cells.Item(99, 12).Value2 = 1878.5
$ws.Cells.Item(99, 13).Value2 = -1879585.8
$ws.Cells.Item(99, 14).Value2 = -4874.5

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Cells.Item(126, 8).Value2 = 1463482.6
$ws.Cells.Item(126, 9).Value2 = 1881083.8
$ws.Cells.Item(126, 10).Value2 = 1878.5
$ws.Cells.Item(126, 11).Value2 = 5643251.4
$ws.Cells.Item(126, 12).Value2 = 5635.5
$ws.Cells.Item(126, 13).Value2 = -5640781.4
$ws.Cells.Item(126, 14).Value2 = -10575.5

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value2 = 6891.222
$ws.Cells.Item(132, 9).Value2 = 7938.8
$ws.Cells.Item(132, 10).Value2 = 1653.3334
$ws.Cells.Item(132, 11).Value2 = 23816.4
$ws.Cells.Item(132, 12).Value2 = 4960.0002
$ws.Cells.Item(132, 13).Value2 = -21286.4
$ws.Cells.Item(132, 14).Value2 = -10020.0002

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Cells.Item(134, 8).Value2 = 7247471
$ws.Cells.Item(134, 9).Value2 = 9524853
$ws.Cells.Item(134, 10).Value2 = 1254.8182
$ws.Cells.Item(134, 11).Value2 = 28574559
$ws.Cells.Item(134, 12).Value2 = 3764.4546
$ws.Cells.Item(134, 13).Value2 = -28572024
$ws.Cells.Item(134, 14).Value2 = -8834.454600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 6: Meat-lover's Special / Marmot Steak
$ws.Cells.Item(6, 8).Value2 = 702.1875
$ws.Cells.Item(6, 9).Value2 = 137.88889
$ws.Cells.Item(6, 10).Value2 = 1427.7142
$ws.Cells.Item(6, 11).Value2 = 413.66667
$ws.Cells.Item(6, 12).Value2 = 4283.142599999999
$ws.Cells.Item(6, 13).Value2 = -300.66667
$ws.Cells.Item(6, 14).Value2 = -4509.142599999999

# Row 68: Such a Butter Face / Fermented Butter
$ws.Cells.Item(68, 8).Value2 = 2559.1562
$ws.Cells.Item(68, 10).Value2 = 2622.3547
$ws.Cells.Item(68, 12).Value2 = 7867.0641
$ws.Cells.Item(68, 14).Value2 = -9489.0641

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Cells.Item(71, 8).Value2 = 2559.1562
$ws.Cells.Item(71, 10).Value2 = 2622.3547
$ws.Cells.Item(71, 12).Value2 = 23601.1923
$ws.Cells.Item(71, 14).Value2 = -31713.1923

# Row 107: Slippery Service / Frantoio Oil
$ws.Cells.Item(107, 8).Value2 = 6728.4707
$ws.Cells.Item(107, 9).Value2 = 675.5
$ws.Cells.Item(107, 10).Value2 = 10030.091
$ws.Cells.Item(107, 11).Value2 = 2026.5
$ws.Cells.Item(107, 12).Value2 = 30090.273
$ws.Cells.Item(107, 13).Value2 = -106.5
$ws.Cells.Item(107, 14).Value2 = -33930.273

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value2 = 24391820
$ws.Cells.Item(131, 9).Value2 = 200000620
$ws.Cells.Item(131, 10).Value2 = 1707.3611
$ws.Cells.Item(131, 11).Value2 = 600001860
$ws.Cells.Item(131, 12).Value2 = 5122.0833
$ws.Cells.Item(131, 13).Value2 = -599996820
$ws.Cells.Item(131, 14).Value2 = -15202.0833

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value2 = 2786.7693
$ws.Cells.Item(126, 9).Value2 = 1546
$ws.Cells.Item(126, 10).Value2 = 3562.25
$ws.Cells.Item(126, 11).Value2 = 4638
$ws.Cells.Item(126, 12).Value2 = 10686.75
$ws.Cells.Item(126, 13).Value2 = -2168
$ws.Cells.Item(126, 14).Value2 = -15626.75

# Row 132: On Board for Lar / Lar Ingot
$ws.Cells.Item(132, 8).Value2 = 2660.7407
$ws.Cells.Item(132, 9).Value2 = 2288.6316
$ws.Cells.Item(132, 10).Value2 = 3544.5
$ws.Cells.Item(132, 11).Value2 = 6865.8948
$ws.Cells.Item(132, 12).Value2 = 10633.5
$ws.Cells.Item(132, 13).Value2 = -4335.8948
$ws.Cells.Item(132, 14).Value2 = -15693.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value2 = 2423.4285
$ws.Cells.Item(40, 9).Value2 = 2423.4285
$ws.Cells.Item(40, 11).Value2 = 2423.4285
$ws.Cells.Item(40, 13).Value2 = -2287.4285

$ws = $wb.Worksheets.Item("WVR")
# Row 43: Walk Softly and Carry a Big Halberd / Velveteen Dress Shoes
$ws.Cells.Item(43, 8).Value2 = 3035
$ws.Cells.Item(43, 9).Value2 = 3035
$ws.Cells.Item(43, 11).Value2 = 3035
$ws.Cells.Item(43, 13).Value2 = -2886

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value2 = 8831.076999999999
$ws.Cells.Item(132, 9).Value2 = 10200.777
$ws.Cells.Item(132, 11).Value2 = 30602.331
$ws.Cells.Item(132, 13).Value2 = -28072.331
